$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (in order) that change for each data row, and the new values per row,
# taken from the target OOXML diff.
$colOrder = @('E','G','H','I','J','K','M','N','O','P','Q','R','S','T')

$updates = @{
    2  = @(3, 18.82499266666667, 56.474978, 0.06886869772378311, 0.0688686977237831, 3, 153.5290173333333, 460.587052, 0.3172206968818489, 0.317220696881849, 2890.182625420539, 26011.64362878486, 0.02184657628528388, 0.02184657628528388)
    3  = @(3, 18.82499266666667, 56.474978, 0.06886869772378311, 0.0688686977237831, 3, 168.7997026666667, 506.3991080000001, 0.3487728915577651, 0.3487728915577651, 3177.653164835514, 28598.87848351963, 0.02401953484294151, 0.0240195348429415)
    4  = @(3, 18.82499266666667, 56.474978, 0.06886869772378311, 0.0688686977237831, 3, 68.09032333333333, 204.27097, 0.1406878008722904, 0.1406878008722904, 1281.799837420962, 11536.19853678866, 0.009688985631697556, 0.009688985631697556)
    5  = @(3, 18.82499266666667, 56.474978, 0.06886869772378311, 0.0688686977237831, 3, 93.562673, 280.688019, 0.1933186106880956, 0.1933186106880956, 1761.316633098731, 15851.84969788858, 0.01331360096386017, 0.01331360096386016)
    6  = @(3, 121.8208923333333, 365.462677, 0.4456653109566078, 0.4456653109566078, 3, 153.5290173333333, 460.587052, 0.3172206968818489, 0.3172206968818489, 18703.04189060647, 168327.3770154582, 0.141374260517721, 0.141374260517721)
    7  = @(3, 121.8208923333333, 365.462677, 0.4456653109566078, 0.4456653109566078, 3, 168.7997026666667, 506.3991080000001, 0.3487728915577651, 0.3487728915577651, 20563.33040445469, 185069.9736400922, 0.1554359791693266, 0.1554359791693266)
    8  = @(3, 121.8208923333333, 365.462677, 0.4456653109566078, 0.4456653109566078, 3, 68.09032333333333, 204.27097, 0.1406878008722904, 0.1406878008722904, 8294.823947731855, 74653.41552958669, 0.06269967252355062, 0.06269967252355062)
    9  = @(3, 121.8208923333333, 365.462677, 0.4456653109566078, 0.4456653109566078, 3, 93.562673, 280.688019, 0.1933186106880956, 0.1933186106880956, 11397.88831395187, 102580.9948255669, 0.08615539874600955, 0.08615539874600954)
    10 = @(3, 87.673585, 263.020755, 0.3207419907481189, 0.3207419907481188, 3, 153.5290173333333, 460.587052, 0.3172206968818489, 0.317220696881849, 13460.43935114047, 121143.9541602643, 0.1017459978243898, 0.1017459978243898)
    11 = @(3, 87.673585, 263.020755, 0.3207419907481189, 0.3207419907481188, 3, 168.7997026666667, 506.3991080000001, 0.3487728915577651, 0.3487728915577651, 14799.27507972073, 133193.4757174866, 0.1118661115572154, 0.1118661115572153)
    12 = @(3, 87.673585, 263.020755, 0.3207419907481189, 0.3207419907481188, 3, 68.09032333333333, 204.27097, 0.1406878008722904, 0.1406878008722904, 5969.722750442484, 53727.50475398235, 0.04512448532575335, 0.04512448532575335)
    13 = @(3, 87.673585, 263.020755, 0.3207419907481189, 0.3207419907481188, 3, 93.562673, 280.688019, 0.1933186106880956, 0.1933186106880956, 8202.974964092706, 73826.77467683435, 0.06200539604076036, 0.06200539604076036)
    14 = @(3, 45.02666966666666, 135.080009, 0.1647240005714903, 0.1647240005714903, 3, 153.5290173333333, 460.587052, 0.3172206968818489, 0.317220696881849, 6912.900347715939, 62216.10312944346, 0.05225386225445422, 0.05225386225445422)
    15 = @(3, 45.02666966666666, 135.080009, 0.1647240005714903, 0.1647240005714903, 3, 168.7997026666667, 506.3991080000001, 0.3487728915577651, 0.3487728915577651, 7600.488451803553, 68404.39606623197, 0.05745126598828161, 0.0574512659882816)
    16 = @(3, 45.02666966666666, 135.080009, 0.1647240005714903, 0.1647240005714903, 3, 68.09032333333333, 204.27097, 0.1406878008722904, 0.1406878008722904, 3065.880496226525, 27592.92446603873, 0.02317465739128887, 0.02317465739128887)
    17 = @(3, 45.02666966666666, 135.080009, 0.1647240005714903, 0.1647240005714903, 3, 93.562673, 280.688019, 0.1933186106880956, 0.1933186106880956, 4212.815570301353, 37915.34013271217, 0.03184421493746557, 0.03184421493746557)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    for ($i = 0; $i -lt $colOrder.Length; $i++) {
        $col = $colOrder[$i]
        $ws.Range("$col$row").Value = $vals[$i]
    }
}
